$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F9").Value = "use restrictions"
$ws.Range("F16").Value = "93_referral_statement"
$ws.Range("F17").Value = "93_referral_statement"
$ws.Range("F30").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F32").Value = "ppe"
$ws.Range("F34").Value = "ppe"
$ws.Range("F35").Value = "application instructions || env warning - water"
$ws.Range("F42").Value = "application instructions"
$ws.Range("F43").Value = "application instructions"
$ws.Range("F44").Value = "application instructions"
$ws.Range("F45").Value = "135_product_information"
$ws.Range("F46").Value = "mixing || application instructions"
$ws.Range("F47").Value = "use restrictions"
$ws.Range("F163").Value = "154_pesticide_storage"
